$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.177.88'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.560.17'
$ws.Range("E3").Value = '  +1.72%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.99'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.73'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.558.68'
$ws.Range("E7").Value = '  +1.74%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.489'
$ws.Range("E9").Value = '  +2.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.136'
$ws.Range("E10").Value = '  +0.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.90'
$ws.Range("E11").Value = '  -3.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.413'
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.160.54'
$ws.Range("E13").Value = '  +1.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000206'
$ws.Range("E14").Value = '  +1.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.02'
$ws.Range("E15").Value = '  -1.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.558.08'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.239.48'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.116'
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.39'
$ws.Range("E19").Value = '  +6.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("E20").Value = '  +0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.74'
$ws.Range("E21").Value = '  -0.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '430.33'
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.614'
$ws.Range("E23").Value = '  +3.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.02'
$ws.Range("E24").Value = '  +1.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.702.18'
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000117'
$ws.Range("E27").Value = '  +0.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.50'
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.07'
$ws.Range("E29").Value = '  -2.13%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.85'
$ws.Range("E30").Value = '  -1.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.51'
$ws.Range("E32").Value = '  +0.98%  '
$ws.Range("B33").Value = 'RenzoRestakedETH'
$ws.Range("C33").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.554.33'
$ws.Range("E33").Value = '  +1.68%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.44'
$ws.Range("E34").Value = '  -2.00%  '
$ws.Range("E35").Value = '  -7.78%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.85'
$ws.Range("E36").Value = '  +1.86%  '
$ws.Range("B37").Value = 'USDe'
$ws.Range("C37").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.73'
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.56'
$ws.Range("E39").Value = '  -0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '174.00'
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0843'
$ws.Range("E41").Value = '  -2.33%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.895'
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.18'
$ws.Range("E43").Value = '  +0.70%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.96'
$ws.Range("E44").Value = '  +3.44%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.08'
$ws.Range("E45").Value = '  +1.61%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.20'
$ws.Range("E47").Value = '  -2.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.95'
$ws.Range("E48").Value = '  -4.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.38'
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.13'
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.06'
$ws.Range("E51").Value = '  +4.84%  '
